# Image_Insert_SQL_Generator.xlsx -- "Update photo insert excel sheet"
#
# 1. Rename the existing sheet "Sheet1" -> "20201003".
# 2. Add a new sheet "20201004" right after it, populated with 25 more
#    photo-insert SQL rows (IDs 101-125) for shop da04f5c9-ffb0-11ea-ba65-065a10bcba76.
# 3. Update the selection/active-cell state on both sheets to match what was
#    left selected when the file was saved.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the first (and, so far, only) sheet -------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "20201003"

# --- 2. Insert the new day's sheet right after it ---------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "20201004"

# Header row, matching sheet1's A/B/C headers (D has no header on sheet1 either).
$ws2.Range("A1").Value = "ID"
$ws2.Range("B1").Value = "Shop ID"
$ws2.Range("C1").Value = "SQL"

$shopId = "da04f5c9-ffb0-11ea-ba65-065a10bcba76"

# Data rows: 25 photos (ids 101-125), all already marked DONE.
for ($i = 0; $i -lt 25; $i++) {
    $row = 2 + $i
    $id = 101 + $i
    $ws2.Range("A$row").Value = $id
    $ws2.Range("B$row").Value = $shopId
    $ws2.Range("D$row").Value = "DONE"
}

# C2 gets the "master" formula; C3:C26 reuse it (Excel records this as one
# shared-formula group anchored at C3, same shape as the existing groups on
# sheet1).
$formulaRow2 = '=_xlfn.CONCAT("INSERT INTO photos(restaurant_id, name, type) VALUES(UuidToBin(''", B2, "''), LPAD(", A2, ", 7, ''0''), ''dish''", ");")'
$formulaRow3 = '=_xlfn.CONCAT("INSERT INTO photos(restaurant_id, name, type) VALUES(UuidToBin(''", B3, "''), LPAD(", A3, ", 7, ''0''), ''dish''", ");")'
$ws2.Range("C2").Formula = $formulaRow2
$ws2.Range("C3:C26").Formula = $formulaRow3

# --- 3. Restore the saved selection on each sheet ---------------------------
# Sheet1 was left with A1:D2 selected (no particular active cell), sheet2
# (the active tab) was left with B12 selected.
$ws1.Range("A1:D2").Select()
$ws2.Activate()
$ws2.Range("B12").Select()
